$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new data row (row 40) with the same shape/format as the
# existing data rows (2-39): a Damasco price record for Macroferia
# Regional de Talca.
$newRow = 40

$ws.Cells.Item($newRow, 1).Value2 = 5
$ws.Cells.Item($newRow, 2).Value2 = "Macroferia Regional de Talca"
$ws.Cells.Item($newRow, 3).Value2 = "Maule"

$ws.Cells.Item($newRow, 4).Value2 = 44890
$ws.Cells.Item($newRow, 4).NumberFormat = $ws.Cells.Item($newRow - 1, 4).NumberFormat

$ws.Cells.Item($newRow, 5).Value2 = 7
$ws.Cells.Item($newRow, 6).Value2 = "Fruta"
$ws.Cells.Item($newRow, 7).Value2 = 100103
$ws.Cells.Item($newRow, 8).Value2 = "Frutos de hueso (carozo)"
$ws.Cells.Item($newRow, 9).Value2 = 100103003
$ws.Cells.Item($newRow, 10).Value2 = "Damasco"
$ws.Cells.Item($newRow, 11).Value2 = "Castle Brite"
$ws.Cells.Item($newRow, 12).Value2 = "Primera"
$ws.Cells.Item($newRow, 13).Value2 = 180
$ws.Cells.Item($newRow, 14).Value2 = 20000
$ws.Cells.Item($newRow, 15).Value2 = 20000
$ws.Cells.Item($newRow, 16).Value2 = 20000
$ws.Cells.Item($newRow, 17).Value2 = "$/caja 12 kilos granel"
$ws.Cells.Item($newRow, 18).Value2 = "Región de O'Higgins"
$ws.Cells.Item($newRow, 19).Value2 = 20000
$ws.Cells.Item($newRow, 20).Value2 = 1
